$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "56.822.40"
$ws.Range("E2").Value = "  +4.49%  "
$ws.Range("D3").Value = "3.247.59"
$ws.Range("E3").Value = "  +2.17%  "
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "395.84"
$ws.Range("E5").Value = "  -1.61%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "108.52"
$ws.Range("E6").Value = "  -0.23%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.581"
$ws.Range("E7").Value = "  +5.66%  "
$ws.Range("D8").Value = "3.244.87"
$ws.Range("E8").Value = "  +2.27%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.999"
$ws.Range("E9").Value = "  -0.04%  "
$ws.Range("E10").Value = "  +0.80%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "39.15"
$ws.Range("E11").Value = "  +0.68%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0971"
$ws.Range("E12").Value = "  +10.19%  "
$ws.Range("E13").Value = "  +1.78%  "
$ws.Range("D14").Value = "3.764.91"
$ws.Range("E14").Value = "  +2.47%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "8.29"
$ws.Range("E15").Value = "  +3.37%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "19.09"
$ws.Range("E16").Value = "  +0.10%  "
$ws.Range("D17").Value = "3.219.44"
$ws.Range("E17").Value = "  +1.13%  "
$ws.Range("E18").Value = "  -3.35%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "10.67"
$ws.Range("E19").Value = "  +1.84%  "
$ws.Range("D20").Value = "56.723.66"
$ws.Range("E20").Value = "  +4.31%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "3.34"
$ws.Range("E21").Value = "  +0.83%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.0000108"
$ws.Range("E22").Value = "  +9.21%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "12.87"
$ws.Range("E23").Value = "  -0.03%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "294.95"
$ws.Range("E24").Value = "  +7.58%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "74.10"
$ws.Range("E25").Value = "  +3.00%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "3.16"
$ws.Range("E26").Value = "  -3.97%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "28.15"
$ws.Range("E27").Value = "  +1.64%  "
$ws.Range("E28").Value = "  +0.89%  "
$ws.Range("E29").Value = "  -4.86%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "7.22"
$ws.Range("E30").Value = "  -2.43%  "
$ws.Range("E31").Value = "  -1.76%  "
$ws.Range("E32").Value = "  +0.05%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "11.31"
$ws.Range("E33").Value = "  +2.10%  "
$ws.Range("E34").Value = "  -4.24%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "39.52"
$ws.Range("E35").Value = "  +6.82%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.0482"
$ws.Range("E36").Value = "  -3.28%  "
$ws.Range("E37").Value = "  +2.21%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "51.49"
$ws.Range("E38").Value = "  +1.32%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.999"
$ws.Range("E39").Value = "  -0.10%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.47"
$ws.Range("E40").Value = "  -4.83%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.88"
$ws.Range("E41").Value = "  +1.90%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "135.63"
$ws.Range("E42").Value = "  +4.07%  "
$ws.Range("E43").Value = "  +3.72%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.88"
$ws.Range("E44").Value = "  -2.96%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "17.01"
$ws.Range("E45").Value = "  -1.45%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.93"
$ws.Range("E46").Value = "  -5.58%  "
$ws.Range("E47").Value = "  -4.24%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "22.10"
$ws.Range("E48").Value = "  -1.47%  "
$ws.Range("E49").Value = "  +3.27%  "
$ws.Range("D50").Value = "2.150.77"
$ws.Range("E50").Value = "  +2.80%  "
$ws.Range("E51").Value = "  -7.23%  "

# Reset number format on cells we forced to text, to avoid leaving
# stray formatting on these General-style cells.
$ws.Range("D5").Style = "Normal"
$ws.Range("D6").Style = "Normal"
$ws.Range("D7").Style = "Normal"
$ws.Range("D9").Style = "Normal"
$ws.Range("D11").Style = "Normal"
$ws.Range("D12").Style = "Normal"
$ws.Range("D15").Style = "Normal"
$ws.Range("D16").Style = "Normal"
$ws.Range("D19").Style = "Normal"
$ws.Range("D21").Style = "Normal"
$ws.Range("D22").Style = "Normal"
$ws.Range("D23").Style = "Normal"
$ws.Range("D24").Style = "Normal"
$ws.Range("D25").Style = "Normal"
$ws.Range("D26").Style = "Normal"
$ws.Range("D27").Style = "Normal"
$ws.Range("D30").Style = "Normal"
$ws.Range("D33").Style = "Normal"
$ws.Range("D35").Style = "Normal"
$ws.Range("D36").Style = "Normal"
$ws.Range("D38").Style = "Normal"
$ws.Range("D39").Style = "Normal"
$ws.Range("D40").Style = "Normal"
$ws.Range("D41").Style = "Normal"
$ws.Range("D42").Style = "Normal"
$ws.Range("D44").Style = "Normal"
$ws.Range("D45").Style = "Normal"
$ws.Range("D46").Style = "Normal"
$ws.Range("D48").Style = "Normal"
